# edit.ps1 - apply the code_smells_element2.docx changes:
#  1. OLEObject/v:shape id + size + ObjectID refresh (re-embedded OLE object)
#  2. "Duplicated code:" paragraph heading -> "Speculative generality" + ":" (2 runs)
#  3. Paragraph body text update: "basically repeated" -> longer phrase,
#     split into 3 runs matching the authored diff.

$d = $word.ActiveDocument

function Find-ParagraphIndexByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# --- 1. OLE object (width/height + shape id + ObjectID) -------------------
# The embedded OLE package is the lone empty paragraph holding a w:object
# field; locate it as the paragraph whose Range contains the document's
# single OLE field.
$oleParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Fields.Count -gt 0) {
        $oleParaIndex = $i
        break
    }
}
if ($oleParaIndex -eq -1) { $oleParaIndex = 7 }

$oleXml = '<w:p w14:paraId="70529D30" w14:textId="36D30125" w:rsidR="009B2757" w:rsidRDefault="00A54893" w:rsidP="009B2757"><w:pPr><w:pStyle w:val="Listacommarcas"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="720"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r w:rsidRPr="00E57372"><w:rPr><w:u w:val="single"/></w:rPr><w:object w:dxaOrig="1740" w:dyaOrig="816" w14:anchorId="7325663C"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1025" type="#_x0000_t75" style="width:87.75pt;height:40.5pt" o:ole=""><v:imagedata r:id="rId9" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Package" ShapeID="_x0000_i1025" DrawAspect="Content" ObjectID="_1700117283" r:id="rId10"/></w:object></w:r></w:p>'
$oleRange = $d.Paragraphs.Item($oleParaIndex).Range
$oleRange.InsertXML($oleXml)
# InsertXML loses the paragraph's direct w:ind left-indent (720 twips =
# 36pt) when the paragraph uses the "Listacommarcas" list style; restore it.
$d.Paragraphs.Item($oleParaIndex).LeftIndent = 36

# --- 2. "Duplicated code:" -> "Speculative generality" + ":" --------------
$dupIndex = Find-ParagraphIndexByText $d "Duplicated code:"
if ($dupIndex -eq -1) { $dupIndex = Find-ParagraphIndexByText $d "Speculative generality" }
if ($dupIndex -ne -1) {
    $dupXml = '<w:p w14:paraId="79465979" w14:textId="113BE6C4" w:rsidR="00184568" w:rsidRDefault="00184568" w:rsidP="00184568"><w:pPr><w:pStyle w:val="Ttulo1"/></w:pPr><w:r><w:t>Speculative generality</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>'
    $d.Paragraphs.Item($dupIndex).Range.InsertXML($dupXml)
}

# --- 3. Expand "basically repeated" into the longer description ----------
$longIndex = Find-ParagraphIndexByText $d "UndoRedoEvent"
if ($longIndex -ne -1) {
    $longXml = '<w:p w14:paraId="4210A998" w14:textId="24DA6608" w:rsidR="00184568" w:rsidRPr="00184568" w:rsidRDefault="00184568" w:rsidP="00184568"><w:r><w:t xml:space="preserve">There is a class inside the logic.undo package called UndoRedoEvent. This is a subclass of the UndoChangeEvent, within the same class. However, this subclass does not have any data or functionalities other than its superclasses’s. So, this class is </w:t></w:r><w:r><w:t>being created perhaps with the intention to add other functionalities in the future but does not have them at the moment</w:t></w:r><w:r><w:t>. A simple solution for this is to delete the subclass and use the superclass instead.</w:t></w:r></w:p>'
    $d.Paragraphs.Item($longIndex).Range.InsertXML($longXml)
}
